$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy header style (bold/centered/bordered) from existing header cell H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$data = @{
    2  = @(8, 9)
    3  = @(9, 9)
    4  = @(7, 7)
    5  = @(10, 10)
    6  = @(9, 9)
    7  = @(9, 9)
    8  = @(9, 9)
    9  = @(7, 8)
    10 = @(7, 7)
    11 = @(8, 9)
    12 = @(9, 9)
    13 = @(5, 6)
    14 = @(9, 9)
    15 = @(8, 8)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
